$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "ProjectParentID" column (column I) entirely; cells to the
# right (J, K, L) shift one column to the left, and the shared string
# "ProjectParentID" is dropped from the workbook since it's no longer used.
$ws.Range("I1:I2").EntireColumn.Delete()
